$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the wide "notes" column from E to F (the new sheet has 6 columns,
#    A..F, instead of the original 5, A..E), matching the column width
#    definitions in the target file. Using Insert() relocates the <col>
#    element verbatim (preserving exact width), then we delete a harmless,
#    unused far-away column to correct the trailing column-range bookkeeping.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()
$ws.Columns("Z:Z").Delete()

# ---------------------------------------------------------------------------
# 2. Build style "templates" in a scratch area far away from the real data,
#    before any of the real target cells are touched.
# ---------------------------------------------------------------------------

# Template A: the pre-existing "date-like, quote-prefixed text" style
# (numFmtId 14 + quotePrefix) -- grab it from B1 before B1 gets overwritten.
$ws.Range("B1").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)

# Template B: "wrap text, vertical centered" style (already used elsewhere
# in the sheet, but build independently so we do not depend on touching the
# original cell before it gets overwritten).
$ws.Range("ZZ2").VerticalAlignment = -4108
$ws.Range("ZZ2").WrapText = $true
$ws.Range("ZZ2").Value = "w"

# Template C: new "plain quote-prefixed text" style (numFmtId 0 + quotePrefix)
$ws.Range("ZZ3").VerticalAlignment = -4108
$ws.Range("ZZ3").Value = "'x"

$dateStyle = $ws.Range("ZZ1")
$wrapStyle = $ws.Range("ZZ2")
$plainStyle = $ws.Range("ZZ3")

# ---------------------------------------------------------------------------
# 3. Write the six rows of activity-delivery data.
# ---------------------------------------------------------------------------
$dates = @("22/10/2021", "25/10/2021", "26/10/2021", "27/10/2021", "28/10/2021", "29/10/2021")

for ($i = 0; $i -lt $dates.Length; $i++) {
  $r = $i + 1

  $wrapStyle.Copy()
  $ws.Range("A$r").PasteSpecial(-4122)
  $ws.Range("F$r").PasteSpecial(-4122)

  $dateStyle.Copy()
  $ws.Range("B$r").PasteSpecial(-4122)

  $plainStyle.Copy()
  $ws.Range("C$r").PasteSpecial(-4122)
  $ws.Range("D$r").PasteSpecial(-4122)
  $ws.Range("E$r").PasteSpecial(-4122)

  $ws.Range("A$r").Value = "Entrega de actividades"
  $ws.Range("B$r").Value = "'" + $dates[$i]
  $ws.Range("C$r").Value = "'09"
  $ws.Range("D$r").Value = "'20"
  $ws.Range("E$r").Value = "'00"
  $ws.Range("F$r").Value = "Entrega de conocimientos tras salida del equipo de System Test"

  $ws.Rows.Item($r).RowHeight = 30
}

# ---------------------------------------------------------------------------
# 4. Remove the scratch template cells.
# ---------------------------------------------------------------------------
$ws.Range("ZZ1:ZZ3").Clear()

# ---------------------------------------------------------------------------
# 5. Selection matches the end-user cursor position in the target file.
# ---------------------------------------------------------------------------
$ws.Range("B6").Select()
